# Applies updated Betfair Back/Lay odds for 2026-01-26 to Sheet1.
# Each assignment below mirrors one changed cell from the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.32
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 3.4
$ws.Range("K2").Value = 3.95
$ws.Range("L2").Value = 1.31
$ws.Range("P2").Value = 2.28
$ws.Range("V2").Value = 1.42
$ws.Range("W2").Value = 1.76

# Row 4
$ws.Range("F4").Value = 2.06
$ws.Range("G4").Value = 2.24
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 4.8
$ws.Range("K4").Value = 3.5
$ws.Range("Q4").Value = 2.3

# Row 5
$ws.Range("F5").Value = 3.9
$ws.Range("G5").Value = 4.6
$ws.Range("H5").Value = 1.93
$ws.Range("I5").Value = 2.3
$ws.Range("K5").Value = 4.5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 3.8
$ws.Range("P5").Value = 1.94
$ws.Range("Q5").Value = 1.87
$ws.Range("R5").Value = 1.38
$ws.Range("S5").Value = 3.15
$ws.Range("T5").Value = 1.74
$ws.Range("U5").Value = 2.08
$ws.Range("V5").Value = 1.76
$ws.Range("W5").Value = 1.27

# Row 8
$ws.Range("Q8").Value = 1.41

# Row 10
$ws.Range("F10").Value = 3.35
$ws.Range("H10").Value = 2.1
$ws.Range("P10").Value = 1.7

# Row 11
$ws.Range("F11").Value = 1.82
$ws.Range("J11").Value = 2.88
$ws.Range("P11").Value = 1.97

# Row 12
$ws.Range("F12").Value = 1.9
$ws.Range("G12").Value = 2.14
$ws.Range("H12").Value = 3.4
$ws.Range("I12").Value = 4.3
$ws.Range("J12").Value = 4.2
$ws.Range("K12").Value = 5.6
$ws.Range("P12").Value = 3.15
$ws.Range("Q12").Value = 1.33

# Row 13
$ws.Range("F13").Value = 2.36
$ws.Range("G13").Value = 2.62
$ws.Range("H13").Value = 2.58
$ws.Range("I13").Value = 2.96
$ws.Range("J13").Value = 3.95
$ws.Range("K13").Value = 4.7
$ws.Range("P13").Value = 2.62
$ws.Range("Q13").Value = 1.41
$ws.Range("R13").Value = 1.67
$ws.Range("S13").Value = 1.97
$ws.Range("T13").Value = 1.4
$ws.Range("U13").Value = 2.8
$ws.Range("X13").Value = 38
$ws.Range("Y13").Value = 24
$ws.Range("Z13").Value = 29
$ws.Range("AA13").Value = 48
$ws.Range("AB13").Value = 22
$ws.Range("AC13").Value = 13
$ws.Range("AD13").Value = 16
$ws.Range("AE13").Value = 29
$ws.Range("AF13").Value = 26
$ws.Range("AG13").Value = 15
$ws.Range("AH13").Value = 16.5
$ws.Range("AI13").Value = 34
$ws.Range("AJ13").Value = 40
$ws.Range("AK13").Value = 26
$ws.Range("AL13").Value = 30
$ws.Range("AM13").Value = 55
$ws.Range("AN13").Value = 12.5
$ws.Range("AO13").Value = 15

# Row 14
$ws.Range("G14").Value = 2.04
$ws.Range("H14").Value = 3.85
$ws.Range("I14").Value = 5.1
$ws.Range("J14").Value = 4
$ws.Range("P14").Value = 2.46

# Row 16
$ws.Range("F16").Value = 1.43
$ws.Range("G16").Value = 1.57
$ws.Range("H16").Value = 6.2
$ws.Range("I16").Value = 8.4
$ws.Range("K16").Value = 6.6
$ws.Range("P16").Value = 2.88

# Row 20
$ws.Range("F20").Value = 1.24
$ws.Range("G20").Value = 2.3
$ws.Range("I20").Value = 5.2
$ws.Range("J20").Value = 2.7

# Row 21
$ws.Range("H21").Value = 2.94
$ws.Range("J21").Value = 2.98

# Row 22
$ws.Range("F22").Value = 3.05
$ws.Range("H22").Value = 2.88
$ws.Range("O22").Value = 1.53
$ws.Range("U22").Value = 1.86
$ws.Range("AF22").Value = 23

# Row 23
$ws.Range("F23").Value = 2.5
$ws.Range("G23").Value = 3.1
$ws.Range("H23").Value = 1.47
$ws.Range("I23").Value = 980
$ws.Range("J23").Value = 1.47
$ws.Range("Q23").Value = 3.4

# Row 24
$ws.Range("F24").Value = 2.6
$ws.Range("G24").Value = 2.62
$ws.Range("H24").Value = 3.15
$ws.Range("N24").Value = 3.2
$ws.Range("T24").Value = 1.96
$ws.Range("AO24").Value = 48

# Row 25
$ws.Range("J25").Value = 3.15
$ws.Range("W25").Value = 1.84
$ws.Range("AE25").Value = 95
